$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> updated numeric value (RandomForest imputation re-run values)
$updates = @{
    "C4" = -12.5955
    "D4" = -8.138800000000005
    "C6" = -11.64549999999999
    "C7" = -12.76640000000001
    "C8" = -10.6216
    "D8" = -7.262900000000005
    "D9" = -7.27159999999999
    "D12" = -5.893699999999999
    "C16" = -14.18869999999999
    "D17" = -8.2303
    "D18" = -9.221699999999993
    "D19" = -8.257099999999996
    "C20" = -11.504
    "D20" = -7.935499999999999
    "C21" = -13.4378
    "D26" = -7.4472
    "C28" = -13.1473
    "C29" = -11.5242
    "C30" = -12.7814
    "D31" = -7.079299999999992
    "C32" = -11.65049999999999
    "D39" = -8.074199999999998
    "C40" = -11.3812
    "D40" = -7.866899999999998
    "D41" = -7.806199999999995
    "D42" = -8.212399999999992
    "D43" = -7.298900000000006
    "C46" = -14.64229999999999
    "D47" = -7.687800000000006
    "D48" = -7.120999999999994
    "C51" = -11.9547
    "C52" = -10.644
    "D54" = -7.838499999999998
    "C57" = -14.4335
    "C59" = -13.16260000000001
    "C62" = -14.0828
    "D62" = -9.35479999999999
    "D63" = -6.555199999999996
    "D64" = -6.977899999999996
    "C66" = -11.3673
    "C73" = -11.88010000000001
    "C74" = -11.7216
    "D76" = -7.632000000000001
    "C77" = -12.757
    "D81" = -7.521199999999998
    "D84" = -8.0053
    "D89" = -8.204999999999993
    "C92" = -10.6655
    "D94" = -5.926299999999999
    "C100" = -11.92689999999999
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
